# Update the five "two-digit divided by one-digit" division problems in
# each of the five data rows of the worksheet table. The table has 20
# rows total; data lives in rows 1, 5, 9, 13, 17 (1-based), the other
# rows being blank spacer rows. Each data row has 5 columns.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# New text values, organized as [row][col], matching the document order
# of the data rows/cells (top-to-bottom, left-to-right).
$newValues = @(
    @("20÷7=2, 6", "95÷8=11, 7", "46÷4=11, 2", "81÷7=11, 4", "68÷6=11, 2"),
    @("87÷9=9, 6", "48÷2=24, 0", "55÷5=11, 0", "69÷2=34, 1", "33÷4=8, 1"),
    @("13÷2=6, 1", "53÷4=13, 1", "25÷4=6, 1", "35÷9=3, 8", "54÷5=10, 4"),
    @("18÷5=3, 3", "30÷3=10, 0", "70÷4=17, 2", "20÷3=6, 2", "39÷6=6, 3"),
    @("78÷7=11, 1", "84÷5=16, 4", "13÷6=2, 1", "76÷8=9, 4", "19÷8=2, 3")
)

$dataRows = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $rowIndex = $dataRows[$i]
    $rowValues = $newValues[$i]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $rowValues[$col - 1]
    }
}
